$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.472.69'
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").Value = '1.837.05'
$ws.Range("E3").Value = '  -0.79%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''261.04'
$ws.Range("E5").Value = '  -0.95%  '

$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").Value = '''0.5385'
$ws.Range("E7").Value = '  +2.26%  '

$ws.Range("D8").Value = '''0.2973'
$ws.Range("E8").Value = '  -8.35%  '

$ws.Range("D9").Value = '''0.06933'
$ws.Range("E9").Value = '  +1.85%  '

$ws.Range("D10").Value = '''17.47'
$ws.Range("E10").Value = '  -7.94%  '

$ws.Range("D11").Value = '1.842.66'
$ws.Range("E11").Value = '  -1.28%  '

$ws.Range("D12").Value = '''0.7311'
$ws.Range("E12").Value = '  -6.79%  '

$ws.Range("D13").Value = '''0.07227'
$ws.Range("E13").Value = '  -6.96%  '

$ws.Range("D14").Value = '''89.04'
$ws.Range("E14").Value = '  +0.26%  '

$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").Value = '''13.81'
$ws.Range("E17").Value = '  -1.43%  '

$ws.Range("E18").Value = '  -0.01%  '

$ws.Range("D19").Value = '''0.000007905'
$ws.Range("E19").Value = '  -0.86%  '

$ws.Range("D20").Value = '26.492.97'
$ws.Range("E20").Value = '  -0.53%  '

$ws.Range("D21").Value = '2.079.15'
$ws.Range("E21").Value = '  -0.81%  '

$ws.Range("D22").Value = '''4.582'
$ws.Range("E22").Value = '  -1.30%  '

$ws.Range("D23").Value = '''5.993'
$ws.Range("E23").Value = '  -0.44%  '

$ws.Range("D24").Value = '''9.202'
$ws.Range("E24").Value = '  -3.08%  '

$ws.Range("D25").Value = '''142.47'
$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("D26").Value = '''2.171'
$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("E27").Value = '  +1.64%  '

$ws.Range("D28").Value = '''16.99'
$ws.Range("E28").Value = '  -0.21%  '

$ws.Range("D29").Value = '''111.01'
$ws.Range("E29").Value = '  -0.61%  '

$ws.Range("D30").Value = '''4.228'
$ws.Range("E30").Value = '  +0.65%  '

$ws.Range("D31").Value = '''0.08872'
$ws.Range("E31").Value = '  +1.68%  '

$ws.Range("E32").Value = '  -1.93%  '

$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("D34").Value = '''2.918'
$ws.Range("E34").Value = '  +1.54%  '

$ws.Range("D35").Value = '''0.7238'
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").Value = '''1.130'
$ws.Range("E36").Value = '  -0.31%  '

$ws.Range("D37").Value = '''3.094'
$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("E38").Value = '  +0.90%  '

$ws.Range("D39").Value = '''0.01708'
$ws.Range("E39").Value = '  -4.68%  '

$ws.Range("D40").Value = '''0.4690'
$ws.Range("E40").Value = '  -3.75%  '

$ws.Range("D41").Value = '''0.9028'
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").Value = '''107.15'
$ws.Range("E42").Value = '  -3.74%  '

$ws.Range("D43").Value = '''5.869'
$ws.Range("E43").Value = '  -1.63%  '

$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''9.075'
$ws.Range("E46").Value = '  +0.51%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '''0.1246'
$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("D48").Value = '''0.4064'
$ws.Range("E48").Value = '  -3.31%  '

$ws.Range("D49").Value = '''34.72'
$ws.Range("E49").Value = '  -1.22%  '

$ws.Range("D50").Value = '''0.05760'
$ws.Range("E50").Value = '  -2.05%  '

$ws.Range("D51").Value = '''0.8908'
$ws.Range("E51").Value = '  -0.02%  '
